$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '327.51'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-1.28%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '44.29'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-0.90%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.294'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-4.58%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08361'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '1.97%'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-5.10%'
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '4.404'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-0.43%'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9701'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-0.71%'
$ws.Range("B9").Value = 'BTSEToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.510'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-4.59%'
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1134'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '0.90%'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1898'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-0.37%'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09668'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-3.78%'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.04603'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-2.03%'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.1060'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.24%'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001300'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '3.04%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005806'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-1.58%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.389'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '1.19%'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '0.18%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.517'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-16.75%'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '0.38%'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '0.26%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.04151'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '1.05%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.001233'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-5.04%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.004418'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.55%'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '1.80%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0002982'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02714'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '-2.29%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05610'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-2.17%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.007859'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '3.20%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1412'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-0.74%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.007318'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-2.75%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002051'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008661'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '4.67%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3516'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '1.36%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000750'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.22%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.003493'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-2.06%'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '40.45%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002101'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.22%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002001'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.22%'
